$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - replace existing data
$ws.Range("A2").Value = "ATHARVA  KADAM"
$ws.Range("B2").Value = "atharva.kadam@ltimindtree.com"
$ws.Range("C2").Value = 1.5
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = "2025-01-02 | 05:06:22 PM"
$ws.Range("F2").Value = "2025-01-02 | 05:04:35 PM"
$ws.Range("G2").Value = "1.78 mins"
$ws.Range("H2").Value = "AI analysis could not be generated due to an error."
$ws.Range("I2").Formula = "="""""

# Row 3
$ws.Range("A3").Value = "Maithri Jajala"
$ws.Range("B3").Value = "maithri.jajala@ltimindtree.com"
$ws.Range("C3").Value = 56.400000000000034
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = "2025-01-09 | 05:01:24 PM"
$ws.Range("F3").Value = "2025-01-09 | 05:00:06 PM"
$ws.Range("G3").Value = "1.30 mins"
$ws.Range("H3").Value = "AI analysis could not be generated due to an error."
$ws.Range("I3").Formula = "="""""

# Row 4
$ws.Range("A4").Value = "Eshwari Kankatte"
$ws.Range("B4").Value = "eshwari@ltimindtree.com"
$ws.Range("C4").Value = 25.2
$ws.Range("D4").Value = 60
$ws.Range("E4").Value = "2025-01-17 | 04:51:43 PM"
$ws.Range("F4").Value = "2025-01-17 | 04:44:54 PM"
$ws.Range("G4").Value = "The difference is more than 5 minutes or not recorded on submission of test. Check manually for Latest Code"
$ws.Range("H4").Value = "AI analysis could not be generated due to an error."
$ws.Range("I4").Value = "https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX19knTC4r7kNG%2F4hDeqFDVN3WwCAuH3MGMHDGumdCpCGc9kxehhTrdANDpBVn7wq8HI4EMz2QzQECDK3DOTwF62a0grY5%2BMPrKAlXJe25F0nXMukxHz%2FmI2HEdj3E6og73X1ftheU%2BXyNg%3D%3D"

# Row 5
$ws.Range("A5").Value = "Saloni Kharat"
$ws.Range("B5").Value = "saloni.kharat@ltimindtree.com"
$ws.Range("C5").Value = 25.2
$ws.Range("D5").Value = 60
$ws.Range("E5").Value = "2025-01-17 | 04:57:40 PM"
$ws.Range("F5").Value = "2025-01-17 | 04:52:59 PM"
$ws.Range("G5").Value = "4.68 mins"
$ws.Range("H5").Value = "AI analysis could not be generated due to an error."
$ws.Range("I5").Formula = "="""""

# Row 6
$ws.Range("A6").Value = "Rohit Patel"
$ws.Range("B6").Value = "rohit.patel@ltimindtree.com"
$ws.Range("C6").Value = 1.2
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = "2024-12-27 | 04:21:20 PM"
$ws.Range("F6").Value = "2024-12-27 | 04:20:53 PM"
$ws.Range("G6").Value = "0.45 mins"
$ws.Range("H6").Value = "AI analysis could not be generated due to an error."
$ws.Range("I6").Formula = "="""""

# Row 7 - note: no SonarAddedTime (F7) for this row
$ws.Range("A7").Value = "Saloni Smriti"
$ws.Range("B7").Value = "saloni.smriti@ltimindtree.com"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = "2024-12-27 | 04:20:10 PM"
$ws.Range("G7").Value = "The difference is more than 5 minutes or not recorded on submission of test. Check manually for Latest Code"
$ws.Range("H7").Value = "No solution is fetched"
$ws.Range("I7").Value = "https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2FVK09YSDwpB2lqU8jC2EdIsC4F5yR%2F4hUn1Ww%2FF0YGAxg7EbRUiLQf2SvTiogw7YI2E6WSWG4X3qHLQDFxOXuF31TriOzf46SrB%2F5HPkwSjq8t6XnB1iQ7wYdcMqnamJN0f9WCwN0Q4Q%3D%3D"
